$wb = $excel.ActiveWorkbook

# sheet1 (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 68
$ws.Range("F4").Value = 961
$ws.Range("G4").Value = 68.8
$ws.Range("F5").Value = 1262
$ws.Range("F6").Value = 1742
$ws.Range("F9").Value = 2598
$ws.Range("F10").Value = 740
$ws.Range("F11").Value = 575
$ws.Range("F12").Value = 571
$ws.Range("F13").Value = 34
$ws.Range("F16").Value = 307
$ws.Range("F18").Value = 2129
$ws.Range("F19").Value = 1233
$ws.Range("F20").Value = 709
$ws.Range("F22").Value = 2622
$ws.Range("F26").Value = 520
$ws.Range("F28").Value = 419
$ws.Range("F36").Value = 341
$ws.Range("F37").Value = 4586
$ws.Range("F38").Value = 154

# sheet2 (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 4202
$ws.Range("F12").Value = 15
$ws.Range("F13").Value = 316
$ws.Range("F14").Value = 326
$ws.Range("F27").Value = 257
$ws.Range("F29").Value = 270
$ws.Range("F38").Value = 478
$ws.Range("F39").Value = 4

# sheet3 (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F7").Value = 132
$ws.Range("F8").Value = 204

# sheet4 (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F6").Value = 68
$ws.Range("F7").Value = 961
$ws.Range("G7").Value = 68.8
$ws.Range("F8").Value = 1262
$ws.Range("F9").Value = 1742
$ws.Range("F12").Value = 132
$ws.Range("F16").Value = 2598
$ws.Range("F17").Value = 740
$ws.Range("F18").Value = 575
$ws.Range("F19").Value = 571
$ws.Range("F20").Value = 34
$ws.Range("F24").Value = 307
$ws.Range("F25").Value = 326
$ws.Range("F27").Value = 2129
$ws.Range("F28").Value = 1233
$ws.Range("F29").Value = 709
$ws.Range("F32").Value = 2622
$ws.Range("F38").Value = 204
$ws.Range("F40").Value = 419
$ws.Range("F41").Value = 419
$ws.Range("F42").Value = 257
$ws.Range("F47").Value = 341
$ws.Range("F48").Value = 4586
$ws.Range("F49").Value = 154
$ws.Range("F50").Value = 478
